# Turn off electrification for building biomass
# On the RBFF sheet, the fuel-shifting matrix has rows = "To type" and
# columns = "From type". The "biomass" column (G) previously sent 100%
# of its fuel to "electricity" (row 2) and 0% to itself (row 7 - biomass).
# This change reverses that: biomass no longer shifts to electricity,
# and instead stays as biomass.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RBFF")

# From biomass -> To electricity: 1 -> 0
$ws.Range("G2").Value = 0

# From biomass -> To biomass: 0 -> 1
$ws.Range("G7").Value = 1

# Reflect the last active selection recorded on the sheet in the saved file
$ws.Range("G8").Select()

# Restore "About" as the active/visible sheet, as it was before the edit
$wb.Worksheets.Item("About").Select()
